$wb = $excel.ActiveWorkbook

# 1. Reorder the worksheet tabs so "review_info" comes before "hotel_info".
$hotelSheet = $wb.Worksheets.Item("hotel_info")
$reviewSheet = $wb.Worksheets.Item("review_info")
$reviewSheet.Move($hotelSheet)

# 2. Insert a new "State" column into hotel_info, between Hotel_Name and City.
$hotelWs = $wb.Worksheets.Item("hotel_info")
$hotelWs.Columns.Item(3).Insert()
$hotelWs.Cells.Item(1, 3).Value = "State"
$hotelWs.Cells.Item(2, 3).Value = "Louisiana"
